{"js": "// Locate the paragraph that currently reads\n// \"        \u041f\u0430\u0440\u0430\u043b\u043b\u0435\u043b\u044c\u043d\u043e\u0441\u0442\u044c \u0432\u0440\u043e\u0434\u0435 \u043d\u0435\u043f\u043b\u043e\u0445\u043e \u0440\u0430\u0431\u043e\u0442\u0430\u0435\u0442.\" (several runs) and\n// replace it with a single run containing \"????\", marking the paragraph\n// mark's run properties with English (US) language, matching the\n// \"save version architecture with Permeatility\" edit.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find((p) => p.text.indexOf(\"\u041f\u0430\u0440\u0430\u043b\u043b\u0435\u043b\") !== -1);\nif (!target) {\n  throw new Error(\"Could not find the target paragraph containing '\u041f\u0430\u0440\u0430\u043b\u043b\u0435\u043b'.\");\n}\n\n// Clear out the existing runs (the multi-run sentence) first.\nconst whole = target.getRange();\nwhole.insertText(\"\", \"Replace\");\nawait context.sync();\n\n// Stamp the (now empty) paragraph's range with English (US) \u2014 with no\n// runs left, this lands on the paragraph mark's run properties\n// (w:pPr/w:rPr/w:lang), not on an individual run.\nconst langRange = target.getRange();\nlangRange.languageId = \"en-US\";\nawait context.sync();\n\n// Insert the replacement text right at the paragraph start, ahead of the\n// \"_GoBack\" bookmark that lives at the end of the paragraph, producing a\n// single clean run with no direct formatting of its own.\nconst start = target.getRange(\"Start\");\nstart.insertText(\"????\", \"Before\");\nawait context.sync();\n", "ps1": "# Locate the paragraph that currently reads\n# \"        \u041f\u0430\u0440\u0430\u043b\u043b\u0435\u043b\u044c\u043d\u043e\u0441\u0442\u044c \u0432\u0440\u043e\u0434\u0435 \u043d\u0435\u043f\u043b\u043e\u0445\u043e \u0440\u0430\u0431\u043e\u0442\u0430\u0435\u0442.\" (several runs) and\n# replace it with a single run containing \"????\", marking the paragraph\n# mark's run properties with English (US) language, matching the\n# \"save version architecture with Permeatility\" edit.\n$d = $word.ActiveDocument\n\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -like \"*\u041f\u0430\u0440\u0430\u043b\u043b\u0435\u043b*\") {\n        $targetIndex = $i\n    }\n}\nif ($targetIndex -eq -1) {\n    throw \"Could not find the target paragraph containing '\u041f\u0430\u0440\u0430\u043b\u043b\u0435\u043b'.\"\n}\n\n$target = $d.Paragraphs($targetIndex)\n\n# Clear out the existing runs (the multi-run sentence) first. Going through\n# $d.Range(start, end) here (rather than $target.Range directly) makes sure\n# the whole paragraph span is replaced, not just its first run.\n$clearRange = $d.Range($target.Range.Start, $target.Range.End)\n$clearRange.Text = \"\"\n\n# Stamp the (now empty) paragraph's range with English (US) \u2014 with no runs\n# left, this lands on the paragraph mark's run properties\n# (w:pPr/w:rPr/w:lang) instead of an individual run.\n$target = $d.Paragraphs($targetIndex)\n$target.Range.LanguageID = \"en-US\"\n\n# Insert the replacement text right at the paragraph start, ahead of the\n# \"_GoBack\" bookmark that lives at the end of the paragraph, producing a\n# single clean run with no direct formatting of its own.\n$target = $d.Paragraphs($targetIndex)\n$startRange = $d.Range($target.Range.Start, $target.Range.Start)\n$startRange.InsertBefore(\"????\")\n"}
